$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the table with the new Guinea master-data layout ---
# Old layout:  A=lang_code B=code C=name D=descr E=is_active
# New layout:  A=code B=name C=descr D=lang_code E=is_active
#              F=cr_by G=cr_dtimes H=upd_by I=upd_dtimes J=is_deleted K=del_dtimes

# Header row
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "lang_code"
$ws.Range("E1").Value = "is_active"
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# Row 2 - FNR / Empreintes digitales
$ws.Range("A2").Value = "FNR"
$ws.Range("B2").Value = "Empreintes digitales"
$ws.Range("C2").Value = "Empreintes digitales du demandeur"
$ws.Range("D2").Value = "fra"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = "superadmin"
$ws.Range("G2").Value = 45079.576834988424
$ws.Range("G2").NumberFormat = "mm:ss.0"
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = "NULL"

# Row 3 - IRS / Iris
$ws.Range("A3").Value = "IRS"
$ws.Range("B3").Value = "Iris"
$ws.Range("C3").Value = "Iris du demandeur"
$ws.Range("D3").Value = "fra"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = "superadmin"
$ws.Range("G3").Value = 45079.576834988424
$ws.Range("G3").NumberFormat = "mm:ss.0"
$ws.Range("H3").Value = "NULL"
$ws.Range("I3").Value = "NULL"
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = "NULL"

# Row 4 - PHT / Photo
$ws.Range("A4").Value = "PHT"
$ws.Range("B4").Value = "Photo"
$ws.Range("C4").Value = "Photo du visage du demandeur"
$ws.Range("D4").Value = "fra"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = "superadmin"
$ws.Range("G4").Value = 45079.576834988424
$ws.Range("G4").NumberFormat = "mm:ss.0"
$ws.Range("H4").Value = "NULL"
$ws.Range("I4").Value = "NULL"
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = "NULL"

# Selected cell moved to F11 in the saved view
[void]$ws.Range("F11").Select()
